$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 91, shifting the existing rows 91-98 down to 92-99.
$ws.Rows.Item(91).Insert()

# Populate the new row 91 with the new weekly price record.
$ws.Cells.Item(91, 1).Value = 8
$ws.Cells.Item(91, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = "2022-07-06"
$ws.Cells.Item(91, 5).Value = 4
$ws.Cells.Item(91, 6).Value = 100112052
$ws.Cells.Item(91, 7).Value = "Albahaca"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 1000
$ws.Cells.Item(91, 11).Value = 3500
$ws.Cells.Item(91, 12).Value = 4000
$ws.Cells.Item(91, 13).Value = 3750
$ws.Cells.Item(91, 14).Value = "`$/paquete"
$ws.Cells.Item(91, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(91, 16).Value = 3750
$ws.Cells.Item(91, 17).Value = 1
$ws.Cells.Item(91, 18).Value = "Hortaliza"
